$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 24.05516451910286
$ws.Cells.Item(2, 5).Value = 25.15244483947754
$ws.Cells.Item(2, 6).Value = 25.53108509993121
$ws.Cells.Item(2, 7).Value = 23.52970760379803
$ws.Cells.Item(2, 8).Value = 4766000000
$ws.Cells.Item(2, 9).Value = "INTC"

$ws.Cells.Item(3, 4).Value = 23.73098190872396
$ws.Cells.Item(3, 5).Value = 22.53236961364746
$ws.Cells.Item(3, 6).Value = 24.06565982080904
$ws.Cells.Item(3, 7).Value = 21.49720384241088
$ws.Cells.Item(3, 8).Value = 4766000000
$ws.Cells.Item(3, 9).Value = "INTC"

$ws.Cells.Item(4, 4).Value = 23.70838143474026
$ws.Cells.Item(4, 5).Value = 26.57285118103028
$ws.Cells.Item(4, 6).Value = 27.49104925146849
$ws.Cells.Item(4, 7).Value = 23.11979200690428
$ws.Cells.Item(4, 8).Value = 4766000000
$ws.Cells.Item(4, 9).Value = "INTC"

$ws.Cells.Item(5, 4).Value = 26.77584172232866
$ws.Cells.Item(5, 5).Value = 24.51554298400879
$ws.Cells.Item(5, 6).Value = 26.87858038148215
$ws.Cells.Item(5, 7).Value = 23.08507342200963
$ws.Cells.Item(5, 8).Value = 4766000000
$ws.Cells.Item(5, 9).Value = "INTC"

$ws.Cells.Item(6, 4).Value = 25.72791577405147
$ws.Cells.Item(6, 5).Value = 24.14134788513184
$ws.Cells.Item(6, 6).Value = 25.88737045984817
$ws.Cells.Item(6, 7).Value = 23.99783897204932
$ws.Cells.Item(6, 8).Value = 4766000000
$ws.Cells.Item(6, 9).Value = "INTC"

$ws.Cells.Item(7, 4).Value = 26.24768370565484
$ws.Cells.Item(7, 5).Value = 28.03291320800781
$ws.Cells.Item(7, 6).Value = 28.89336094330404
$ws.Cells.Item(7, 7).Value = 26.03860423766414
$ws.Cells.Item(7, 8).Value = 4766000000
$ws.Cells.Item(7, 9).Value = "INTC"

$ws.Cells.Item(8, 4).Value = 30.53841147815202
$ws.Cells.Item(8, 5).Value = 28.25349998474121
$ws.Cells.Item(8, 6).Value = 31.08128218691229
$ws.Cells.Item(8, 7).Value = 28.12385972592597
$ws.Cells.Item(8, 8).Value = 4766000000
$ws.Cells.Item(8, 9).Value = "INTC"

$ws.Cells.Item(9, 4).Value = 29.88794089172134
$ws.Cells.Item(9, 5).Value = 30.05938148498535
$ws.Cells.Item(9, 6).Value = 31.39009371566189
$ws.Cells.Item(9, 7).Value = 29.54505659092279
$ws.Cells.Item(9, 8).Value = 4766000000
$ws.Cells.Item(9, 9).Value = "INTC"

$ws.Cells.Item(10, 4).Value = 29.75597545178443
$ws.Cells.Item(10, 5).Value = 29.72308921813965
$ws.Cells.Item(10, 6).Value = 30.8330788003511
$ws.Cells.Item(10, 7).Value = 28.95020470218228
$ws.Cells.Item(10, 8).Value = 4766000000
$ws.Cells.Item(10, 9).Value = "INTC"

$ws.Cells.Item(11, 4).Value = 27.75740442580016
$ws.Cells.Item(11, 5).Value = 29.38093757629395
$ws.Cells.Item(11, 6).Value = 29.7039865564943
$ws.Cells.Item(11, 7).Value = 27.52547247280213
$ws.Cells.Item(11, 8).Value = 4766000000
$ws.Cells.Item(11, 9).Value = "INTC"

$ws.Cells.Item(12, 4).Value = 31.81305736445479
$ws.Cells.Item(12, 5).Value = 37.96369552612305
$ws.Cells.Item(12, 6).Value = 38.22240408831388
$ws.Cells.Item(12, 7).Value = 31.7796777710977
$ws.Cells.Item(12, 8).Value = 4766000000
$ws.Cells.Item(12, 9).Value = "INTC"

$ws.Cells.Item(13, 4).Value = 38.93582612807229
$ws.Cells.Item(13, 5).Value = 40.41333770751953
$ws.Cells.Item(13, 6).Value = 42.68837114284218
$ws.Cells.Item(13, 7).Value = 35.62820977469026
$ws.Cells.Item(13, 8).Value = 4766000000
$ws.Cells.Item(13, 9).Value = "INTC"

$ws.Cells.Item(14, 4).Value = 43.68793204707458
$ws.Cells.Item(14, 5).Value = 43.62876892089844
$ws.Cells.Item(14, 6).Value = 47.15321790789084
$ws.Cells.Item(14, 7).Value = 39.96909269477364
$ws.Cells.Item(14, 8).Value = 4766000000
$ws.Cells.Item(14, 9).Value = "INTC"

$ws.Cells.Item(15, 4).Value = 41.68739959083602
$ws.Cells.Item(15, 5).Value = 40.88833236694336
$ws.Cells.Item(15, 6).Value = 45.30869341160787
$ws.Cells.Item(15, 7).Value = 39.46871818086888
$ws.Cells.Item(15, 8).Value = 4766000000
$ws.Cells.Item(15, 9).Value = "INTC"

$ws.Cells.Item(16, 4).Value = 40.01662578662575
$ws.Cells.Item(16, 5).Value = 40.09359741210938
$ws.Cells.Item(16, 6).Value = 42.52247582878741
$ws.Cells.Item(16, 7).Value = 36.22791749472168
$ws.Cells.Item(16, 8).Value = 4766000000
$ws.Cells.Item(16, 9).Value = "INTC"

$ws.Cells.Item(17, 4).Value = 39.55415811381391
$ws.Cells.Item(17, 5).Value = 40.5524787902832
$ws.Cells.Item(17, 6).Value = 43.06549494464316
$ws.Cells.Item(17, 7).Value = 38.20298280039454
$ws.Cells.Item(17, 8).Value = 4766000000
$ws.Cells.Item(17, 9).Value = "INTC"

$ws.Cells.Item(18, 4).Value = 47.06259577350248
$ws.Cells.Item(18, 5).Value = 44.20454406738281
$ws.Cells.Item(18, 6).Value = 51.60949726627182
$ws.Cells.Item(18, 7).Value = 44.0573093486784
$ws.Cells.Item(18, 8).Value = 4766000000
$ws.Cells.Item(18, 9).Value = "INTC"

$ws.Cells.Item(19, 4).Value = 42.95034468047347
$ws.Cells.Item(19, 5).Value = 44.04828262329102
$ws.Cells.Item(19, 6).Value = 46.61885570550808
$ws.Cells.Item(19, 7).Value = 41.01587895011813
$ws.Cells.Item(19, 8).Value = 4766000000
$ws.Cells.Item(19, 9).Value = "INTC"

$ws.Cells.Item(20, 4).Value = 45.59139636543285
$ws.Cells.Item(20, 5).Value = 49.59171676635742
$ws.Cells.Item(20, 6).Value = 50.21457654769927
$ws.Cells.Item(20, 7).Value = 42.57360704235445
$ws.Cells.Item(20, 8).Value = 4766000000
$ws.Cells.Item(20, 9).Value = "INTC"

$ws.Cells.Item(21, 4).Value = 53.13721770987009
$ws.Cells.Item(21, 5).Value = 56.39213562011719
$ws.Cells.Item(21, 6).Value = 61.12014875792387
$ws.Cells.Item(21, 7).Value = 51.62001856014938
$ws.Cells.Item(21, 8).Value = 4766000000
$ws.Cells.Item(21, 9).Value = "INTC"

$ws.Cells.Item(22, 4).Value = 46.53788618324653
$ws.Cells.Item(22, 5).Value = 53.16842651367188
$ws.Cells.Item(22, 6).Value = 55.07426510996061
$ws.Cells.Item(22, 7).Value = 45.40324927837772
$ws.Cells.Item(22, 8).Value = 4766000000
$ws.Cells.Item(22, 9).Value = "INTC"

$ws.Cells.Item(23, 4).Value = 53.40635776074765
$ws.Cells.Item(23, 5).Value = 42.54858016967773
$ws.Cells.Item(23, 6).Value = 55.20707329069204
$ws.Cells.Item(23, 7).Value = 41.87108488740854
$ws.Cells.Item(23, 8).Value = 4766000000
$ws.Cells.Item(23, 9).Value = "INTC"

$ws.Cells.Item(24, 4).Value = 47.02885664323092
$ws.Cells.Item(24, 5).Value = 39.74117660522461
$ws.Cells.Item(24, 6).Value = 50.46626928464421
$ws.Cells.Item(24, 7).Value = 39.13985509909842
$ws.Cells.Item(24, 8).Value = 4766000000
$ws.Cells.Item(24, 9).Value = "INTC"

$ws.Cells.Item(25, 4).Value = 45.10181718142243
$ws.Cells.Item(25, 5).Value = 50.18243789672852
$ws.Cells.Item(25, 6).Value = 57.81241287728636
$ws.Cells.Item(25, 7).Value = 44.59556527084455
$ws.Cells.Item(25, 8).Value = 4766000000
$ws.Cells.Item(25, 9).Value = "INTC"

$ws.Cells.Item(26, 4).Value = 58.86365876910687
$ws.Cells.Item(26, 5).Value = 52.32426071166992
$ws.Cells.Item(26, 6).Value = 62.29251834532009
$ws.Cells.Item(26, 7).Value = 51.9058857328104
$ws.Cells.Item(26, 8).Value = 4766000000
$ws.Cells.Item(26, 9).Value = "INTC"

$ws.Cells.Item(27, 4).Value = 51.36537598092723
$ws.Cells.Item(27, 5).Value = 49.15995025634766
$ws.Cells.Item(27, 6).Value = 52.5824764060846
$ws.Cells.Item(27, 7).Value = 47.86963899824732
$ws.Cells.Item(27, 8).Value = 4766000000
$ws.Cells.Item(27, 9).Value = "INTC"

$ws.Cells.Item(28, 4).Value = 49.41493188918744
$ws.Cells.Item(28, 5).Value = 45.13199615478516
$ws.Cells.Item(28, 6).Value = 51.64389817603289
$ws.Cells.Item(28, 7).Value = 44.09119607598761
$ws.Cells.Item(28, 8).Value = 4766000000
$ws.Cells.Item(28, 9).Value = "INTC"

$ws.Cells.Item(29, 4).Value = 47.90363929400464
$ws.Cells.Item(29, 5).Value = 45.27890777587891
$ws.Cells.Item(29, 6).Value = 52.19780602793664
$ws.Cells.Item(29, 7).Value = 42.94169210534397
$ws.Cells.Item(29, 8).Value = 4766000000
$ws.Cells.Item(29, 9).Value = "INTC"

$ws.Cells.Item(30, 4).Value = 46.56771153462336
$ws.Cells.Item(30, 5).Value = 40.73623275756836
$ws.Cells.Item(30, 6).Value = 46.6331284617015
$ws.Cells.Item(30, 7).Value = 40.65212477061715
$ws.Cells.Item(30, 8).Value = 4766000000
$ws.Cells.Item(30, 9).Value = "INTC"

$ws.Cells.Item(31, 4).Value = 34.93576430551267
$ws.Cells.Item(31, 5).Value = 34.20106887817383
$ws.Cells.Item(31, 6).Value = 38.36434775727346
$ws.Cells.Item(31, 7).Value = 33.19321616877606
$ws.Cells.Item(31, 8).Value = 4766000000
$ws.Cells.Item(31, 9).Value = "INTC"

$ws.Cells.Item(32, 4).Value = 25.05109445904402
$ws.Cells.Item(32, 5).Value = 27.04909324645996
$ws.Cells.Item(32, 6).Value = 27.80072028154468
$ws.Cells.Item(32, 7).Value = 23.39561026795782
$ws.Cells.Item(32, 8).Value = 4766000000
$ws.Cells.Item(32, 9).Value = "INTC"

$ws.Cells.Item(33, 4).Value = 26.08371312247101
$ws.Cells.Item(33, 5).Value = 27.25049018859864
$ws.Cells.Item(33, 6).Value = 29.40082918909817
$ws.Cells.Item(33, 7).Value = 25.42800577929437
$ws.Cells.Item(33, 8).Value = 4766000000
$ws.Cells.Item(33, 9).Value = "INTC"

$ws.Cells.Item(34, 4).Value = 32.15034333521978
$ws.Cells.Item(34, 5).Value = 30.3154125213623
$ws.Cells.Item(34, 6).Value = 33.02876954243178
$ws.Cells.Item(34, 7).Value = 27.82654323192737
$ws.Cells.Item(34, 8).Value = 4766000000
$ws.Cells.Item(34, 9).Value = "INTC"

$ws.Cells.Item(35, 4).Value = 32.87001578170341
$ws.Cells.Item(35, 5).Value = 35.05546951293945
$ws.Cells.Item(35, 6).Value = 36.44710221863846
$ws.Cells.Item(35, 7).Value = 30.75316158285264
$ws.Cells.Item(35, 8).Value = 4766000000
$ws.Cells.Item(35, 9).Value = "INTC"

$ws.Cells.Item(36, 4).Value = 35.02421911933925
$ws.Cells.Item(36, 5).Value = 35.89957809448242
$ws.Cells.Item(36, 6).Value = 36.60773535614693
$ws.Cells.Item(36, 7).Value = 31.62113672646627
$ws.Cells.Item(36, 8).Value = 4766000000
$ws.Cells.Item(36, 9).Value = "INTC"

$ws.Cells.Item(37, 4).Value = 48.54977989358274
$ws.Cells.Item(37, 5).Value = 42.51066207885742
$ws.Cells.Item(37, 6).Value = 49.63524092963385
$ws.Cells.Item(37, 7).Value = 41.80017630149584
$ws.Cells.Item(37, 8).Value = 4766000000
$ws.Cells.Item(37, 9).Value = "INTC"

$ws.Cells.Item(38, 4).Value = 43.8327771951091
$ws.Cells.Item(38, 5).Value = 30.15544509887696
$ws.Cells.Item(38, 6).Value = 44.94121392180273
$ws.Cells.Item(38, 7).Value = 30.10596202430713
$ws.Cells.Item(38, 8).Value = 4766000000
$ws.Cells.Item(38, 9).Value = "INTC"

$ws.Cells.Item(39, 4).Value = 30.71515496192047
$ws.Cells.Item(39, 5).Value = 30.54622650146484
$ws.Cells.Item(39, 6).Value = 36.9257573384868
$ws.Cells.Item(39, 7).Value = 29.82082859876
$ws.Cells.Item(39, 8).Value = 4766000000
$ws.Cells.Item(39, 9).Value = "INTC"

$ws.Cells.Item(40, 4).Value = 23.45999908447266
$ws.Cells.Item(40, 5).Value = 21.52000045776367
$ws.Cells.Item(40, 6).Value = 23.81999969482422
$ws.Cells.Item(40, 7).Value = 21.46999931335449
$ws.Cells.Item(40, 8).Value = 4766000000
$ws.Cells.Item(40, 9).Value = "INTC"

$ws.Cells.Item(41, 4).Value = 20.22999954223633
$ws.Cells.Item(41, 5).Value = 19.43000030517578
$ws.Cells.Item(41, 6).Value = 22.40999984741211
$ws.Cells.Item(41, 7).Value = 18.72999954223633
$ws.Cells.Item(41, 8).Value = 4766000000
$ws.Cells.Item(41, 9).Value = "INTC"

$ws.Cells.Item(42, 4).Value = 22.55999946594238
$ws.Cells.Item(42, 5).Value = 20.10000038146973
$ws.Cells.Item(42, 6).Value = 23.89999961853028
$ws.Cells.Item(42, 7).Value = 17.67000007629395
$ws.Cells.Item(42, 8).Value = 4766000000
$ws.Cells.Item(42, 9).Value = "INTC"

$ws.Cells.Item(43, 4).Value = 22.32999992370605
$ws.Cells.Item(43, 5).Value = 19.79999923706055
$ws.Cells.Item(43, 6).Value = 24.04000091552734
$ws.Cells.Item(43, 7).Value = 19.65999984741211
$ws.Cells.Item(43, 8).Value = 4766000000
$ws.Cells.Item(43, 9).Value = "INTC"
